$wb = $excel.ActiveWorkbook

# --- Content change -------------------------------------------------
# The shared string "Ready for handoff" becomes "In Translation".
# It is used by:
#   - "Overview" sheet, columns E (zh-cn) and F (de-de), rows 2-4
#   - "zh-cn" sheet, column C (Status), rows 2-4
#   - "de-de" sheet, column C (Status), rows 2-4
$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$wsOverview.Range("E2:F4").Value = "In Translation"
$wsZhCn.Range("C2:C4").Value = "In Translation"
$wsDeDe.Range("C2:C4").Value = "In Translation"

# --- Column width follow-up -----------------------------------------
# With the shorter replacement text in place, the (now too-wide) status
# columns are re-sized to fit the new content, same as the source report
# generator does when it rebuilds the sheet.
$newWidth = 12.5

$wsOverview.Columns.Item(5).ColumnWidth = $newWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newWidth
$wsZhCn.Columns.Item(3).ColumnWidth = $newWidth
$wsDeDe.Columns.Item(3).ColumnWidth = $newWidth
